$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest coinranking.com snapshot: refreshed Price (D) / Volume(1h) (E) for
# every coin row, plus rows 14/15 which swapped rank (Polygon now ranks
# above Polkadot, so their Coin name + Link also need to be updated).
$updates = @(
    @{ Cell = 'D2'; Value = '25.817.99' },
    @{ Cell = 'E2'; Value = '  -2.51%  ' },
    @{ Cell = 'D3'; Value = '1.749.13' },
    @{ Cell = 'E3'; Value = '  -4.77%  ' },
    @{ Cell = 'D4'; Value = '0.9996' },
    @{ Cell = 'E4'; Value = '  -0.10%  ' },
    @{ Cell = 'D5'; Value = '239.32' },
    @{ Cell = 'E5'; Value = '  -8.48%  ' },
    @{ Cell = 'D6'; Value = '0.9995' },
    @{ Cell = 'E6'; Value = '  -0.12%  ' },
    @{ Cell = 'D7'; Value = '0.5085' },
    @{ Cell = 'E7'; Value = '  -5.32%  ' },
    @{ Cell = 'D8'; Value = '42.06' },
    @{ Cell = 'E8'; Value = '  -6.25%  ' },
    @{ Cell = 'D9'; Value = '0.2765' },
    @{ Cell = 'E9'; Value = '  -7.32%  ' },
    @{ Cell = 'D10'; Value = '0.06181' },
    @{ Cell = 'E10'; Value = '  -10.78%  ' },
    @{ Cell = 'D11'; Value = '1.747.22' },
    @{ Cell = 'E11'; Value = '  -4.97%  ' },
    @{ Cell = 'D12'; Value = '0.06948' },
    @{ Cell = 'E12'; Value = '  -3.54%  ' },
    @{ Cell = 'D13'; Value = '15.68' },
    @{ Cell = 'E13'; Value = '  -10.49%  ' },
    @{ Cell = 'B14'; Value = 'Polygon' },
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' },
    @{ Cell = 'D14'; Value = '0.6030' },
    @{ Cell = 'E14'; Value = '  -17.57%  ' },
    @{ Cell = 'B15'; Value = 'Polkadot' },
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' },
    @{ Cell = 'D15'; Value = '4.519' },
    @{ Cell = 'E15'; Value = '  -9.35%  ' },
    @{ Cell = 'D16'; Value = '77.52' },
    @{ Cell = 'E16'; Value = '  -13.02%  ' },
    @{ Cell = 'D17'; Value = '0.9990' },
    @{ Cell = 'E17'; Value = '  -0.20%  ' },
    @{ Cell = 'D18'; Value = '0.9996' },
    @{ Cell = 'E18'; Value = '  -0.08%  ' },
    @{ Cell = 'D19'; Value = '25.815.39' },
    @{ Cell = 'E19'; Value = '  -2.58%  ' },
    @{ Cell = 'D20'; Value = '0.000006939' },
    @{ Cell = 'E20'; Value = '  -12.09%  ' },
    @{ Cell = 'E21'; Value = '  -15.52%  ' },
    @{ Cell = 'D22'; Value = '1.970.00' },
    @{ Cell = 'E22'; Value = '  -5.14%  ' },
    @{ Cell = 'D23'; Value = '4.077' },
    @{ Cell = 'E23'; Value = '  -10.85%  ' },
    @{ Cell = 'D24'; Value = '5.260' },
    @{ Cell = 'E24'; Value = '  -12.12%  ' },
    @{ Cell = 'D25'; Value = '8.204' },
    @{ Cell = 'E25'; Value = '  -10.69%  ' },
    @{ Cell = 'D26'; Value = '137.91' },
    @{ Cell = 'E26'; Value = '  -3.32%  ' },
    @{ Cell = 'E27'; Value = '  -14.30%  ' },
    @{ Cell = 'D28'; Value = '1.824' },
    @{ Cell = 'E28'; Value = '  -15.78%  ' },
    @{ Cell = 'D29'; Value = '15.03' },
    @{ Cell = 'E29'; Value = '  -11.36%  ' },
    @{ Cell = 'D30'; Value = '103.81' },
    @{ Cell = 'E30'; Value = '  -6.40%  ' },
    @{ Cell = 'D31'; Value = '0.08149' },
    @{ Cell = 'E31'; Value = '  -7.85%  ' },
    @{ Cell = 'D32'; Value = '3.714' },
    @{ Cell = 'E32'; Value = '  -12.09%  ' },
    @{ Cell = 'D33'; Value = '3.501' },
    @{ Cell = 'E33'; Value = '  -13.26%  ' },
    @{ Cell = 'D34'; Value = '0.04534' },
    @{ Cell = 'E34'; Value = '  -6.25%  ' },
    @{ Cell = 'D35'; Value = '0.9987' },
    @{ Cell = 'E35'; Value = '  -0.11%  ' },
    @{ Cell = 'D36'; Value = '2.624' },
    @{ Cell = 'E36'; Value = '  -10.53%  ' },
    @{ Cell = 'D37'; Value = '0.9859' },
    @{ Cell = 'E37'; Value = '  -12.76%  ' },
    @{ Cell = 'D38'; Value = '0.6115' },
    @{ Cell = 'E38'; Value = '  -15.36%  ' },
    @{ Cell = 'D39'; Value = '2.682' },
    @{ Cell = 'E39'; Value = '  -13.24%  ' },
    @{ Cell = 'E40'; Value = '  -9.36%  ' },
    @{ Cell = 'D41'; Value = '1.925' },
    @{ Cell = 'E41'; Value = '  -15.64%  ' },
    @{ Cell = 'D42'; Value = '0.9995' },
    @{ Cell = 'E42'; Value = '  -0.06%  ' },
    @{ Cell = 'D43'; Value = '102.82' },
    @{ Cell = 'E43'; Value = '  -4.51%  ' },
    @{ Cell = 'D44'; Value = '0.3865' },
    @{ Cell = 'E44'; Value = '  -17.76%  ' },
    @{ Cell = 'D45'; Value = '0.7453' },
    @{ Cell = 'E45'; Value = '  -17.50%  ' },
    @{ Cell = 'D46'; Value = '4.927' },
    @{ Cell = 'E46'; Value = '  -16.35%  ' },
    @{ Cell = 'E47'; Value = '  -6.42%  ' },
    @{ Cell = 'E48'; Value = '  -10.59%  ' },
    @{ Cell = 'D49'; Value = '5.991' },
    @{ Cell = 'E49'; Value = '  -18.93%  ' },
    @{ Cell = 'D50'; Value = '30.15' },
    @{ Cell = 'E50'; Value = '  -13.28%  ' },
    @{ Cell = 'D51'; Value = '52.61' },
    @{ Cell = 'E51'; Value = '  -12.30%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Cell.StartsWith("D")) {
        # Price text often looks numeric ("0.9990", "25.817.99", ...). Force
        # the cell to text first so Excel doesn't silently coerce it to a
        # number and drop significant trailing zeros, then drop the number
        # format back to the sheet's default (no explicit style) so the
        # cell's formatting matches the rest of the column.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
